$wb = $excel.ActiveWorkbook

# --- Sheet: "Status by Landings (Area)" ---
$wsArea = $wb.Worksheets.Item("Status by Landings (Area)")
$wsArea.Range("C3").Value = 3.699075689035789
$wsArea.Range("C5").Value = 4.517698030842167
$wsArea.Range("C7").Value = 14.93902192521486
$wsArea.Range("C8").Value = 67.50435457160448
$wsArea.Range("C9").Value = 17.55662350318067
$wsArea.Range("C10").Value = 82.44337649681934
$wsArea.Range("C11").Value = 17.55662350318067

# --- Sheet: "Status by Landings (Tier)" ---
$wsTier = $wb.Worksheets.Item("Status by Landings (Tier)")

# Row 4 ("Area")
$wsTier.Range("C4").Value = 3.699075689035789
$wsTier.Range("E4").Value = 4.517698030842167
$wsTier.Range("G4").Value = 14.93902192521486
$wsTier.Range("H4").Value = 67.50435457160448
$wsTier.Range("I4").Value = 17.55662350318067
$wsTier.Range("J4").Value = 82.44337649681934
$wsTier.Range("K4").Value = 17.55662350318067

# Row 5 ("Global")
$wsTier.Range("C5").Value = 3.699075689035789
$wsTier.Range("E5").Value = 4.517698030842167
$wsTier.Range("G5").Value = 14.93902192521486
$wsTier.Range("H5").Value = 67.50435457160448
$wsTier.Range("I5").Value = 17.55662350318067
$wsTier.Range("J5").Value = 82.44337649681934
$wsTier.Range("K5").Value = 17.55662350318067

# --- Sheet: "Comparison by Landings" ---
$wsComp = $wb.Worksheets.Item("Comparison by Landings")
$wsComp.Range("C2").Value = 98.4871602932466
$wsComp.Range("C3").Value = 14.93902192521486
$wsComp.Range("C4").Value = 67.50435457160448
$wsComp.Range("C5").Value = 17.55662350318067
$wsComp.Range("C6").Value = 82.44337649681934
$wsComp.Range("C7").Value = 17.55662350318067
